$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("D6").Value = "[Python - 프로그래머스] 대충 만든 자판"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-%EB%8C%80%EC%B6%A9-%EB%A7%8C%EB%93%A0-%EC%9E%90%ED%8C%90"

# Row 32
$ws.Range("D32").Value = "SOLID-python 원칙 : clean code"
$ws.Range("E32").Value = "https://dodonam.tistory.com/440"

# Row 36
$ws.Range("D36").Value = "Anomaly Detection in Surveillance Video"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/413"

# Row 51
$ws.Range("D51").Value = "[docker] 도커 이미지 빌드시 이미지 태그명 빌드 시간으로 하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/docker-%EB%8F%84%EC%BB%A4-%EC%9D%B4%EB%AF%B8%EC%A7%80-%EB%B9%8C%EB%93%9C%EC%8B%9C-%EC%9D%B4%EB%AF%B8%EC%A7%80-%ED%83%9C%EA%B7%B8%EB%AA%85-%EB%B9%8C%EB%93%9C-%EC%8B%9C%EA%B0%84%EC%9C%BC%EB%A1%9C-%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
